$d = $word.ActiveDocument

$replacements = @(
    @{old="199÷7="; new="756÷6="},
    @{old="514÷5="; new="691÷6="},
    @{old="486÷9="; new="894÷5="},
    @{old="878÷4="; new="620÷4="},
    @{old="275÷2="; new="629÷6="},
    @{old="946÷8="; new="432÷9="},
    @{old="623÷2="; new="610÷5="},
    @{old="377÷7="; new="226÷5="},
    @{old="810÷9="; new="197÷4="},
    @{old="794÷3="; new="638÷6="},
    @{old="337÷8="; new="597÷6="},
    @{old="461÷6="; new="546÷8="},
    @{old="201÷2="; new="325÷3="},
    @{old="195÷3="; new="867÷6="},
    @{old="558÷7="; new="705÷4="},
    @{old="232÷9="; new="302÷5="},
    @{old="898÷4="; new="347÷6="},
    @{old="405÷7="; new="448÷4="},
    @{old="545÷2="; new="931÷8="},
    @{old="786÷8="; new="834÷4="},
    @{old="887÷3="; new="729÷2="},
    @{old="450÷4="; new="351÷3="},
    @{old="822÷7="; new="549÷3="},
    @{old="756÷9="; new="857÷3="},
    @{old="147÷7="; new="769÷2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
